# Tasks.xlsx update: append "present to a team" follow-up step to the
# "Done criteria" column (E) for several tasks, matching the revised
# presentation requirement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value  = "1. ODATA is working 2. present to a team"
$ws.Range("E5").Value  = "1. team review 2. hi-fi mockup uploaded in github in editable format in design folder 3. Each screen uploaded in document folder as jpg 3. present to a team"
$ws.Range("E6").Value  = "1. team review 2. present to a team"
$ws.Range("E7").Value  = "1. team review 2.  uploaded in github in editable format in design folder 3.uploaded in document folder as jpg 4. present to  ateam"
$ws.Range("E9").Value  = "1. team review 2. uploaded to github 3. present to a team"
$ws.Range("E10").Value = "1. team review 2. all figures are uploaded in github in editable format in design folder 3. all figures are uploaded in document folder as jpg 4. present to a team"
$ws.Range("E11").Value = "1. team review 2. uploaded to github 3. present to a team"
$ws.Range("E12").Value = "1. team review 2. uploaded to github 3. present to a team"

# Row heights grow to fit the longer wrapped text in a few rows.
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 75

# Selection moved from B13 to E13.
$ws.Range("E13").Select()
